# Adds a new "Quantidade Vendida" column (D) to the product/price sheet.
# - D1 gets the same header formatting as A1:C1 (bold, centered, bordered).
# - D2:D121 are populated with 0 (quantity sold not yet tracked).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, formatted like the existing header row (copy C1's style
# onto D1 so it matches Produto/Valor Unitário/Quantidade exactly).
$ws.Range("D1").Value = "Quantidade Vendida"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# New data column, default (unstyled) numeric cells, all starting at 0.
$ws.Range("D2:D121").Value = 0
